$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value is a plain decimal number (e.g. "591.30") must be
# forced to Text format first, otherwise Excel auto-converts the string into a
# numeric value (losing the original text representation, e.g. trailing zeros).
$textFormatRows = 5,6,11,14,16,18,20,21,22,24,26,28,31,32,36,37,40,41,44,45,46,47,48,49
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply the updated coin list values scraped on Sun Jun 23 09:08:48 UTC 2024.
# Row 2
$ws.Cells.Item(2, 4).Value = "64.405.84"
$ws.Cells.Item(2, 5).Value = "  +0.20%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "3.505.93"
$ws.Cells.Item(3, 5).Value = "  +0.58%  "
# Row 5
$ws.Cells.Item(5, 4).Value = "591.30"
$ws.Cells.Item(5, 5).Value = "  +0.90%  "
# Row 6
$ws.Cells.Item(6, 4).Value = "134.68"
$ws.Cells.Item(6, 5).Value = "  +0.62%  "
# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.75%  "
# Row 9
$ws.Cells.Item(9, 5).Value = "  +5.98%  "
# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.32%  "
# Row 11
$ws.Cells.Item(11, 4).Value = "0.391"
$ws.Cells.Item(11, 5).Value = "  +4.17%  "
# Row 12
$ws.Cells.Item(12, 4).Value = "4.102.14"
$ws.Cells.Item(12, 5).Value = "  +0.60%  "
# Row 13
$ws.Cells.Item(13, 5).Value = "  +0.83%  "
# Row 14
$ws.Cells.Item(14, 4).Value = "0.0000182"
$ws.Cells.Item(14, 5).Value = "  +1.09%  "
# Row 15
$ws.Cells.Item(15, 4).Value = "3.506.18"
$ws.Cells.Item(15, 5).Value = "  +0.63%  "
# Row 16
$ws.Cells.Item(16, 4).Value = "25.85"
$ws.Cells.Item(16, 5).Value = "  +2.56%  "
# Row 17
$ws.Cells.Item(17, 4).Value = "64.398.72"
$ws.Cells.Item(17, 5).Value = "  +0.10%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "10.06"
$ws.Cells.Item(18, 5).Value = "  +0.78%  "
# Row 19
$ws.Cells.Item(19, 5).Value = "  +1.88%  "
# Row 20
$ws.Cells.Item(20, 4).Value = "13.62"
$ws.Cells.Item(20, 5).Value = "  -0.39%  "
# Row 21
$ws.Cells.Item(21, 4).Value = "392.11"
$ws.Cells.Item(21, 5).Value = "  +1.97%  "
# Row 22
$ws.Cells.Item(22, 4).Value = "0.583"
$ws.Cells.Item(22, 5).Value = "  +3.17%  "
# Row 23
$ws.Cells.Item(23, 4).Value = "3.645.73"
$ws.Cells.Item(23, 5).Value = "  +0.61%  "
# Row 24
$ws.Cells.Item(24, 4).Value = "74.52"
$ws.Cells.Item(24, 5).Value = "  +0.58%  "
# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.09%  "
# Row 26
$ws.Cells.Item(26, 4).Value = "5.67"
$ws.Cells.Item(26, 5).Value = "  -0.43%  "
# Row 27
$ws.Cells.Item(27, 5).Value = "  +4.06%  "
# Row 28
$ws.Cells.Item(28, 4).Value = "1.01"
$ws.Cells.Item(28, 5).Value = "  +0.60%  "
# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.15%  "
# Row 30
$ws.Cells.Item(30, 5).Value = "  +2.33%  "
# Row 31
$ws.Cells.Item(31, 4).Value = "8.22"
$ws.Cells.Item(31, 5).Value = "  +0.31%  "
# Row 32
$ws.Cells.Item(32, 4).Value = "1.47"
$ws.Cells.Item(32, 5).Value = "  -4.83%  "
# Row 33
$ws.Cells.Item(33, 5).Value = "  +7.46%  "
# Row 34
$ws.Cells.Item(34, 4).Value = "3.533.83"
$ws.Cells.Item(34, 5).Value = "  +0.81%  "
# Row 36
$ws.Cells.Item(36, 4).Value = "23.44"
$ws.Cells.Item(36, 5).Value = "  +0.38%  "
# Row 37
$ws.Cells.Item(37, 4).Value = "5.36"
$ws.Cells.Item(37, 5).Value = "  +1.51%  "
# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.08%  "
# Row 39
$ws.Cells.Item(39, 5).Value = "  +2.78%  "
# Row 40
$ws.Cells.Item(40, 4).Value = "165.83"
$ws.Cells.Item(40, 5).Value = "  +2.21%  "
# Row 41
$ws.Cells.Item(41, 4).Value = "0.0791"
$ws.Cells.Item(41, 5).Value = "  +1.79%  "
# Row 42
$ws.Cells.Item(42, 5).Value = "  +1.14%  "
# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.00%  "
# Row 44
$ws.Cells.Item(44, 4).Value = "4.46"
$ws.Cells.Item(44, 5).Value = "  +1.82%  "
# Row 45
$ws.Cells.Item(45, 4).Value = "25.02"
$ws.Cells.Item(45, 5).Value = "  -2.12%  "
# Row 46
$ws.Cells.Item(46, 2).Value = "ONDO"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(46, 4).Value = "1.18"
$ws.Cells.Item(46, 5).Value = "  -0.92%  "
# Row 47
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).Value = "1.66"
$ws.Cells.Item(47, 5).Value = "  +0.89%  "
# Row 48
$ws.Cells.Item(48, 2).Value = "SuiNetwork"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(48, 4).Value = "0.925"
$ws.Cells.Item(48, 5).Value = "  +3.00%  "
# Row 49
$ws.Cells.Item(49, 2).Value = "Cosmos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(49, 4).Value = "6.82"
$ws.Cells.Item(49, 5).Value = "  +1.38%  "
# Row 50
$ws.Cells.Item(50, 4).Value = "2.391.38"
$ws.Cells.Item(50, 5).Value = "  -2.82%  "
# Row 51
$ws.Cells.Item(51, 5).Value = "  +0.87%  "
